$wb = $excel.ActiveWorkbook

# --- "kriteria" sheet: add the CONCATENATE helper column (G) for rows 19-41 ---
$wsKriteria = $wb.Worksheets.Item("kriteria")

# G19 is its own (non-shared) formula
$wsKriteria.Range("G19").Formula = "=CONCATENATE(""'"",F19,""',"")"

# G20:G41 share one formula (Excel collapses a filled-down range into a shared formula)
$wsKriteria.Range("G20:G41").Formula = "=CONCATENATE(""'"",F20,""',"")"

# Widen column G slightly now that it holds data (closest reachable width to the
# recorded best-fit width of 8.42578125)
$wsKriteria.Columns(7).ColumnWidth = 7.59

# --- "rtrw" sheet: selection moves to A2 ---
$wsRtrw = $wb.Worksheets.Item("rtrw")
$wsRtrw.Activate()
$wsRtrw.Range("A2").Select() | Out-Null

# --- "kriteria" sheet: selection moves to G20 (no longer the active tab) ---
$wsKriteria.Activate()
$wsKriteria.Range("G20").Select() | Out-Null

# --- "kumuhRT" sheet: becomes the active tab, selection moves to AG32 ---
$wsKumuhRT = $wb.Worksheets.Item("kumuhRT")
$wsKumuhRT.Activate()
$wsKumuhRT.Range("AG32").Select() | Out-Null
